# Apply edits to the "poland_division-2_2023-2024" sheet:
#  - Re-order several match rows (rows 34/35, 39/40, 41/42 and 56/57 are swapped;
#    rows 49/50/52 are rotated) by exchanging their F:V (match-data) contents,
#    while leaving the A:E (index/tournament/date) columns untouched.
#  - Append a new match row (row 90) with the GKS Jastrzebie vs Ol. Grudziadz data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper: swap the F:V contents of two rows -----------------------------
function Swap-Rows($sheet, $rowA, $rowB) {
    $rangeA = $sheet.Range("F$rowA`:V$rowA")
    $rangeB = $sheet.Range("F$rowB`:V$rowB")
    $valA = $rangeA.Value()
    $valB = $rangeB.Value()
    $rangeA.Value = $valB
    $rangeB.Value = $valA
}

# Simple pairwise swaps
Swap-Rows $ws 34 35
Swap-Rows $ws 39 40
Swap-Rows $ws 41 42
Swap-Rows $ws 56 57

# 3-way rotation among rows 49, 50, 52:
#   new_49 = old_52 ; new_50 = old_49 ; new_52 = old_50
$range49 = $ws.Range("F49:V49")
$range50 = $ws.Range("F50:V50")
$range52 = $ws.Range("F52:V52")

$val49 = $range49.Value()
$val50 = $range50.Value()
$val52 = $range52.Value()

$range49.Value = $val52
$range50.Value = $val49
$range52.Value = $val50

# --- Add new row 90 ----------------------------------------------------------
# Copy the formatting (styles) of row 89 down into row 90 first, then set values.
$ws.Range("A89:V89").Copy()
$ws.Range("A90:V90").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A90").Value = 89
$ws.Range("B90").Value = "poland"
$ws.Range("C90").Value = "division-2"
$ws.Range("D90").Value = "2023-2024"
$ws.Range("E90").Value = 45194.84375
$ws.Range("F90").Value = "GKS Jastrzebie"
$ws.Range("G90").Value = 3
$ws.Range("H90").Value = "Ol. Grudziadz"
$ws.Range("I90").Value = 2
$ws.Range("J90").Value = 2.07
$ws.Range("K90").Value = "24/09/2023 07:43"
$ws.Range("L90").Value = 2.24
$ws.Range("M90").Value = "25/09/2023 18:23"
$ws.Range("N90").Value = 3.18
$ws.Range("O90").Value = "24/09/2023 07:43"
$ws.Range("P90").Value = 3.38
$ws.Range("Q90").Value = "25/09/2023 20:00"
$ws.Range("R90").Value = 3.17
$ws.Range("S90").Value = "24/09/2023 07:43"
$ws.Range("T90").Value = 3.03
$ws.Range("U90").Value = "25/09/2023 20:00"
$ws.Range("V90").Value = "https://www.betexplorer.com/football/poland/division-2/gks-jastrzebie-ol-grudziadz/6epbYdeM/"

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()
